$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the format from G1 (bold header style)
# then overwrite the pasted value with the correct header text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill H2:H5 with value 0 (plain numeric cells, same as other data columns)
$ws.Range("H2:H5").Value = 0
